$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-01 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-02 Thursday", 2) | Out-Null
$d.Content.Find.Execute("982×6=5892", $true, $false, $false, $false, $false, $true, 1, $false, "698×8=5584", 2) | Out-Null
$d.Content.Find.Execute("348×7=2436", $true, $false, $false, $false, $false, $true, 1, $false, "925×4=3700", 2) | Out-Null
$d.Content.Find.Execute("333×7=2331", $true, $false, $false, $false, $false, $true, 1, $false, "299×5=1495", 2) | Out-Null
$d.Content.Find.Execute("485×6=2910", $true, $false, $false, $false, $false, $true, 1, $false, "912×9=8208", 2) | Out-Null
$d.Content.Find.Execute("236×9=2124", $true, $false, $false, $false, $false, $true, 1, $false, "781×7=5467", 2) | Out-Null
$d.Content.Find.Execute("187×3=561", $true, $false, $false, $false, $false, $true, 1, $false, "221×6=1326", 2) | Out-Null
$d.Content.Find.Execute("402×8=3216", $true, $false, $false, $false, $false, $true, 1, $false, "673×6=4038", 2) | Out-Null
$d.Content.Find.Execute("341×7=2387", $true, $false, $false, $false, $false, $true, 1, $false, "485×9=4365", 2) | Out-Null
$d.Content.Find.Execute("117×9=1053", $true, $false, $false, $false, $false, $true, 1, $false, "546×2=1092", 2) | Out-Null
$d.Content.Find.Execute("195×2=390", $true, $false, $false, $false, $false, $true, 1, $false, "519×5=2595", 2) | Out-Null
$d.Content.Find.Execute("482×3=1446", $true, $false, $false, $false, $false, $true, 1, $false, "756×7=5292", 2) | Out-Null
$d.Content.Find.Execute("338×6=2028", $true, $false, $false, $false, $false, $true, 1, $false, "772×2=1544", 2) | Out-Null
$d.Content.Find.Execute("731×9=6579", $true, $false, $false, $false, $false, $true, 1, $false, "143×5=715", 2) | Out-Null
$d.Content.Find.Execute("496×7=3472", $true, $false, $false, $false, $false, $true, 1, $false, "980×5=4900", 2) | Out-Null
$d.Content.Find.Execute("744×9=6696", $true, $false, $false, $false, $false, $true, 1, $false, "406×7=2842", 2) | Out-Null
$d.Content.Find.Execute("413×8=3304", $true, $false, $false, $false, $false, $true, 1, $false, "791×5=3955", 2) | Out-Null
$d.Content.Find.Execute("895×2=1790", $true, $false, $false, $false, $false, $true, 1, $false, "610×6=3660", 2) | Out-Null
$d.Content.Find.Execute("667×7=4669", $true, $false, $false, $false, $false, $true, 1, $false, "879×8=7032", 2) | Out-Null
$d.Content.Find.Execute("413×5=2065", $true, $false, $false, $false, $false, $true, 1, $false, "841×2=1682", 2) | Out-Null
$d.Content.Find.Execute("320×7=2240", $true, $false, $false, $false, $false, $true, 1, $false, "663×9=5967", 2) | Out-Null
$d.Content.Find.Execute("608×7=4256", $true, $false, $false, $false, $false, $true, 1, $false, "194×6=1164", 2) | Out-Null
$d.Content.Find.Execute("105×3=315", $true, $false, $false, $false, $false, $true, 1, $false, "440×6=2640", 2) | Out-Null
$d.Content.Find.Execute("844×5=4220", $true, $false, $false, $false, $false, $true, 1, $false, "217×8=1736", 2) | Out-Null
$d.Content.Find.Execute("219×4=876", $true, $false, $false, $false, $false, $true, 1, $false, "144×9=1296", 2) | Out-Null
$d.Content.Find.Execute("317×6=1902", $true, $false, $false, $false, $false, $true, 1, $false, "807×2=1614", 2) | Out-Null
